# Calculate T6 enzyme activity
# - Add a new (empty) "T5" sheet after "T3".
# - On "T0" and "T3": add PPO / PER columns (I, J) with header style matching
#   the other headers, and replace the single free-text legend cell (A21)
#   with four separate legend lines placed in column K next to the rows
#   they annotate.
# - On "T0" also fill in the PPO ("a") / PER ("x"/"y") observation data for
#   every sample row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Add the new T5 worksheet as the last tab (after T3).
# ---------------------------------------------------------------------
$ws0 = $wb.Worksheets.Item("T0")
$ws3 = $wb.Worksheets.Item("T3")

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws5 = $wb.Worksheets.Add($null, $lastSheet)
$ws5.Name = "T5"

# ---------------------------------------------------------------------
# 2. T0 sheet: new PPO / PER headers + data, legend in column K, drop the
#    old combined legend cell at A21.
# ---------------------------------------------------------------------
$ws0.Range("H1").Copy()
$ws0.Range("I1").PasteSpecial(-4122)   # xlPasteFormats
$ws0.Range("J1").PasteSpecial(-4122)   # xlPasteFormats
$ws0.Application.CutCopyMode = $false

$ws0.Range("I1").Value = "PPO"
$ws0.Range("J1").Value = "PER"

$ws0_data = @{
    "I2"  = "a";  "J2"  = "y"
    "I3"  = "a";  "J3"  = "x"
    "I4"  = "a";  "J4"  = "y"
    "I5"  = "a";  "J5"  = "x"
    "I6"  = "a";  "J6"  = "x"
    "I7"  = "a";  "J7"  = "x"
    "J8"  = "y"
    "J9"  = "y"
    "I10" = "a";  "J10" = "y"
    "I11" = "a";  "J11" = "y"
    "I12" = "a";  "J12" = "y"
    "J13" = "y"
    "I14" = "a";  "J14" = "y"
    "I15" = "a";  "J15" = "y"
    "I16" = "a";  "J16" = "y"
    "J17" = "y"
}
foreach ($addr in $ws0_data.Keys) {
    $ws0.Range($addr).Value = $ws0_data[$addr]
}

$ws0.Range("K2").Value = "x = noisy, indicating low activity or bad data"
$ws0.Range("K3").Value = "o = possible substrate inhibition"
$ws0.Range("K4").Value = "a = salvageable with other errors"
$ws0.Range("K5").Value = "y = generally negative activity, indicating no activity or need to improve methodology"

$ws0.Range("A21").ClearContents()

# ---------------------------------------------------------------------
# 3. T3 sheet: same new headers (no PPO/PER data entered yet) + legend.
# ---------------------------------------------------------------------
$ws3.Range("H1").Copy()
$ws3.Range("I1").PasteSpecial(-4122)   # xlPasteFormats
$ws3.Range("J1").PasteSpecial(-4122)   # xlPasteFormats
$ws3.Application.CutCopyMode = $false

$ws3.Range("I1").Value = "PPO"
$ws3.Range("J1").Value = "PER"

$ws3.Range("K2").Value = "x = noisy, indicating low activity or bad data"
$ws3.Range("K3").Value = "o = possible substrate inhibition"
$ws3.Range("K4").Value = "a = salvageable with other errors"
$ws3.Range("K5").Value = "y = generally negative activity, indicating no activity or need to improve methodology"

$ws3.Range("A21").ClearContents()

# ---------------------------------------------------------------------
# 4. Restore the selections on each sheet, and make T3 the active tab
#    again (matching the original tabSelected state).
# ---------------------------------------------------------------------
$ws0.Activate()
[void]$ws0.Range("K5").Select()

$ws3.Activate()
[void]$ws3.Range("I2").Select()
